# Apply odds updates to Sheet1 as described in the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 1.9
$ws.Range("I4").Value = 4.2
$ws.Range("U4").Value = 9.5
$ws.Range("AA4").Value = 6.5
$ws.Range("AF4").Value = 21
$ws.Range("AG4").Value = 13

# Row 9
$ws.Range("I9").Value = 3.75
$ws.Range("K9").Value = 5.9
$ws.Range("L9").Value = 1.38
$ws.Range("M9").Value = 2.82
$ws.Range("N9").Value = 2.12
$ws.Range("O9").Value = 1.65
$ws.Range("P9").Value = 1.44
$ws.Range("Q9").Value = 2.62
$ws.Range("R9").Value = 1.75
$ws.Range("S9").Value = 1.98
$ws.Range("T9").Value = 6.9
$ws.Range("U9").Value = 11
$ws.Range("W9").Value = 25
$ws.Range("X9").Value = 19
$ws.Range("Y9").Value = 28
$ws.Range("Z9").Value = 5.9

# Row 12
$ws.Range("K12").Value = 10

# Row 14
$ws.Range("G14").Value = 2.15
$ws.Range("I14").Value = 3.1
$ws.Range("AE14").Value = 12

# Row 17
$ws.Range("G17").Value = 3.5
$ws.Range("H17").Value = 3.05
$ws.Range("I17").Value = 2.1
$ws.Range("R17").Value = 1.83
$ws.Range("S17").Value = 1.78
$ws.Range("U17").Value = 18
$ws.Range("V17").Value = 12
$ws.Range("X17").Value = 35
$ws.Range("AE17").Value = 6.5
$ws.Range("AG17").Value = 8.75
$ws.Range("AH17").Value = 19.5
$ws.Range("AI17").Value = 18
$ws.Range("AJ17").Value = 32

# Row 20
$ws.Range("G20").Value = 2.4
$ws.Range("I20").Value = 2.87
$ws.Range("T20").Value = 6.9
$ws.Range("U20").Value = 11
$ws.Range("W20").Value = 25
$ws.Range("X20").Value = 22
$ws.Range("Y20").Value = 37
$ws.Range("AE20").Value = 7.8
$ws.Range("AG20").Value = 10.75
$ws.Range("AH20").Value = 35
$ws.Range("AI20").Value = 27

# Row 21
$ws.Range("G21").Value = 2.8
$ws.Range("I21").Value = 2.8
$ws.Range("V21").Value = 12
$ws.Range("X21").Value = 29
$ws.Range("AA21").Value = 5.5

# Row 24
$ws.Range("AB24").Value = 11

# Row 33
$ws.Range("G33").Value = 1.6
$ws.Range("I33").Value = 4.2
$ws.Range("AE33").Value = 21
$ws.Range("AG33").Value = 15
$ws.Range("AJ33").Value = 26

# Row 42
$ws.Range("G42").Value = 2.05
$ws.Range("I42").Value = 3.3
$ws.Range("L42").Value = 1.22
$ws.Range("M42").Value = 4
$ws.Range("P42").Value = 1.33
$ws.Range("Q42").Value = 3.25
$ws.Range("R42").Value = 1.67
$ws.Range("S42").Value = 2.1
$ws.Range("T42").Value = 8.5
$ws.Range("U42").Value = 10
$ws.Range("W42").Value = 19
$ws.Range("Y42").Value = 23
$ws.Range("Z42").Value = 12
$ws.Range("AB42").Value = 13
$ws.Range("AC42").Value = 41
$ws.Range("AD42").Value = 151
$ws.Range("AE42").Value = 12
$ws.Range("AI42").Value = 26
$ws.Range("AJ42").Value = 29

# Row 43
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 3.6
$ws.Range("I43").Value = 3.5
$ws.Range("U43").Value = 9.5
$ws.Range("W43").Value = 17
$ws.Range("X43").Value = 15
$ws.Range("AA43").Value = 7
$ws.Range("AC43").Value = 51
$ws.Range("AF43").Value = 19
$ws.Range("AI43").Value = 29

# Row 44
$ws.Range("J44").Value = 1.06
$ws.Range("K44").Value = 10
$ws.Range("L44").Value = 1.33
$ws.Range("M44").Value = 3.25
$ws.Range("N44").Value = 2.08
$ws.Range("O44").Value = 1.73

# Row 46
$ws.Range("H46").Value = 6
$ws.Range("I46").Value = 1.27
$ws.Range("N46").Value = 1.22
$ws.Range("O46").Value = 4.2
